$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Forces the cell to stay text even when the value looks numeric
    # (prices like "5.80" or "1.00" would otherwise be auto-converted
    # by Excel into plain numbers and lose their original formatting).
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Update price (D) and volume/1h change (E) columns for rows with refreshed crypto data
Set-TextValue "D2" "60.503.02"
$ws.Range("E2").Value = "  +4.22%  "
Set-TextValue "D3" "2.448.00"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "556.69"
$ws.Range("E5").Value = "  +3.00%  "
Set-TextValue "D6" "139.38"
$ws.Range("E6").Value = "  +2.44%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.59%  "
Set-TextValue "D9" "0.107"
$ws.Range("E9").Value = "  +4.48%  "
Set-TextValue "D10" "5.80"
$ws.Range("E10").Value = "  +4.27%  "
Set-TextValue "D11" "0.362"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  +4.72%  "
Set-TextValue "D14" "2.886.80"
$ws.Range("E14").Value = "  +4.07%  "
Set-TextValue "D15" "60.393.78"
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("E16").Value = "  +5.21%  "
Set-TextValue "D17" "2.428.55"
$ws.Range("E17").Value = "  +2.81%  "
Set-TextValue "D18" "11.51"
$ws.Range("E18").Value = "  +7.44%  "
Set-TextValue "D19" "4.44"
$ws.Range("E19").Value = "  +3.76%  "
Set-TextValue "D20" "336.55"
$ws.Range("E20").Value = "  +1.32%  "
Set-TextValue "D21" "6.90"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("E22").Value = "  +0.08%  "
Set-TextValue "D23" "64.91"
$ws.Range("E24").Value = "  +2.18%  "
Set-TextValue "D25" "8.61"
$ws.Range("E25").Value = "  +1.03%  "
Set-TextValue "D26" "0.998"
$ws.Range("E26").Value = "  -0.51%  "
Set-TextValue "D27" "1.38"
$ws.Range("E27").Value = "  +0.18%  "
Set-TextValue "D28" "0.0₃0799"
$ws.Range("E28").Value = "  +8.39%  "
Set-TextValue "D29" "1.81"
$ws.Range("E29").Value = "  +3.59%  "
Set-TextValue "D30" "6.32"
$ws.Range("E30").Value = "  +2.88%  "
Set-TextValue "D31" "170.95"
$ws.Range("E31").Value = "  -0.88%  "
Set-TextValue "D32" "18.87"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E35").Value = "  +5.80%  "
Set-TextValue "D36" "4.30"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +0.65%  "
Set-TextValue "D39" "40.13"
$ws.Range("E39").Value = "  +2.02%  "
Set-TextValue "D40" "0.420"
$ws.Range("E40").Value = "  +11.11%  "
Set-TextValue "D41" "316.91"
$ws.Range("E41").Value = "  +7.93%  "
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  +6.55%  "

# Rows 42-46 were re-ranked: Filecoin, Aave, Stellar, InjectiveProtocol, Hedera now occupy these rows
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.75"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "144.56"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D44" "0.0966"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "19.93"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0526"
$ws.Range("E46").Value = "  +4.53%  "
